# First Commit-General opening Cash
# Adds a new "GeneralOpening" worksheet (after the last existing sheet)
# with a header row + one data row, mirroring the layout/style of the
# existing "AccOpn_LoanAdva_PersnlLoanWeek" sheet, and updates that
# sheet's selection now that it is no longer the active tab.

$wb = $excel.ActiveWorkbook

# The sheet that is currently last / active in the workbook.
$lastSheet = $wb.Worksheets.Item("AccOpn_LoanAdva_PersnlLoanWeek")

# Add the new sheet right after it, so it becomes the last tab.
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "GeneralOpening"

# --- Columns A-D (header + data), values already present elsewhere -----
$newSheet.Range("A1").Value = "TestScenario"
$newSheet.Range("A2").Value = "Personal Loan Weekly"
$newSheet.Range("B1").Value = "Run"
$newSheet.Range("B2").Value = "Yes"
$newSheet.Range("C1").Value = "pcRegFormName"
$newSheet.Range("C2").Value = "qwerty"
$newSheet.Range("D1").Value = "pcRegFormPcName"
$newSheet.Range("D2").Value = "zxcvb"

# --- Columns E-H, entered column by column so new shared-string entries
#     land in the same order as the authored workbook (Name, AAA,
#     Openaccount, Remarks, Testing, Amount).
$newSheet.Range("E1").Value = "Name"
$newSheet.Range("E2").Value = "AAA"
$newSheet.Range("F1").Value = "Openaccount"
$newSheet.Range("F2").Value = 500000
$newSheet.Range("G1").Value = "Remarks"
$newSheet.Range("G2").Value = "Testing"
$newSheet.Range("H1").Value = "Amount"
$newSheet.Range("H2").Value = 2000

# --- Formatting: header = bold white-on-blue, wrapped, top/left aligned -
$headerRange = $newSheet.Range("A1:H1")
$headerRange.Font.Bold = $true
$headerRange.Font.Color = 16777215
$headerRange.Interior.Color = 9851952
$headerRange.HorizontalAlignment = -4131
$headerRange.VerticalAlignment = -4160
$headerRange.WrapText = $true

# --- Formatting: data row text cells = wrapped, top/left aligned --------
$newSheet.Range("A2:E2").HorizontalAlignment = -4131
$newSheet.Range("A2:E2").VerticalAlignment = -4160
$newSheet.Range("A2:E2").WrapText = $true

$newSheet.Range("G2").HorizontalAlignment = -4131
$newSheet.Range("G2").VerticalAlignment = -4160
$newSheet.Range("G2").WrapText = $true

# Row heights to match the source sheet's look (both rows tall enough to
# show wrapped header/body text).
$newSheet.Rows(1).RowHeight = 45
$newSheet.Rows(2).RowHeight = 45

# Match page orientation used elsewhere in the workbook.
$newSheet.PageSetup.Orientation = 1

# The new sheet's selection lands on H1.
$newSheet.Range("H1").Select() | Out-Null

# The previously-active sheet keeps a plain range selection now that
# focus has moved to the new tab.
$lastSheet.Range("A1:D2").Select() | Out-Null
$newSheet.Select() | Out-Null
